# Add two new columns, I (I0) and J (IF), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style from H1 onto I1/J1, set header text.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for I and J columns, rows 2-32.
$values = @(
    @(9, 9),
    @(7, 7),
    @(7, 8),
    @(9, 9),
    @(7, 7),
    @(9, 9),
    @(8, 9),
    @(8, 9),
    @(8, 8),
    @(6, 6),
    @(5, 5),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(5, 6),
    @(1, 3),
    @(1, 4),
    @(8, 8),
    @(4, 6),
    @(2, 3),
    @(7, 8),
    @(6, 8),
    @(7, 9),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(10, 10),
    @(4, 6),
    @(5, 5),
    @(7, 9),
    @(9, 9)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
